# ============================================================================
# Lab 3 / rectifier data.xlsx -- apply commit "Report done." edits
# ============================================================================
$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet1 "Halfwave"
# ----------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Swap columns A and B (data + column width/bestFit) using a cut/insert so
# the stored column widths come out exactly swapped (14.5703125 / 14.7109375).
$ws1.Columns.Item(1).Cut()
$ws1.Columns.Item(3).Insert()

# Re-assert exact cell values/headers (the cut/insert above can introduce
# floating point drift, and we still need to fix up the data regardless).
$ws1.Range("A1").Value = "Output Current"
$ws1.Range("B1").Value = "Output Voltage"
$ws1.Range("C1").Value = "P-P Ripple"
$ws1.Range("D1").Value = "Peak Current (Secondary)"

$ws1.Cells.Item(2,1).Value = 0.053
$ws1.Cells.Item(2,2).Value = 14.06
$ws1.Cells.Item(2,3).Value = 0
$ws1.Cells.Item(2,4).Value = 0.085

$ws1.Cells.Item(3,1).Value = 0.501
$ws1.Cells.Item(3,2).Value = 12.3
$ws1.Cells.Item(3,3).Value = 2.05
$ws1.Cells.Item(3,4).Value = 3.5

$ws1.Cells.Item(4,1).Value = 1.02
$ws1.Cells.Item(4,2).Value = 11.5
$ws1.Cells.Item(4,3).Value = 4.2
$ws1.Cells.Item(4,4).Value = 6.5

$ws1.Cells.Item(5,1).Value = 1.588
$ws1.Cells.Item(5,2).Value = 10.6
$ws1.Cells.Item(5,3).Value = 6.5
$ws1.Cells.Item(5,4).Value = 9.4

$ws1.Cells.Item(6,1).Value = 2.07
$ws1.Cells.Item(6,2).Value = 9.86
$ws1.Cells.Item(6,3).Value = 8
$ws1.Cells.Item(6,4).Value = 11.6

# Clear the leftover B8:B11 values (old column-A tail, now shifted to B by
# the cut/insert) and drop the now-empty rows.
$ws1.Range("B8:B11").ClearContents()
$ws1.Rows("8:11").Delete()

# New column E: "R" header + condition/resistance values.
$ws1.Range("E1").Value = "R"
$ws1.Range("E2").Value = "Inf"
$ws1.Cells.Item(3,5).Value = 25
$ws1.Cells.Item(4,5).Value = 12
$ws1.Cells.Item(5,5).Value = 7
$ws1.Cells.Item(6,5).Value = 4

$ws1.PageSetup.Orientation = 1
$ws1.Range("E1:E6").Select()

# ----------------------------------------------------------------------
# Sheet2 "Fullwave"
# ----------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(3,1).Value = 12.009
$ws2.Cells.Item(3,2).Value = 0.566
$ws2.Cells.Item(3,3).Value = 0.5
$ws2.Cells.Item(3,4).Value = 1.8

$ws2.Cells.Item(4,1).Value = 11.68
$ws2.Cells.Item(4,2).Value = 1.017
$ws2.Cells.Item(4,3).Value = 0.9
$ws2.Cells.Item(4,4).Value = 1.82

$ws2.Cells.Item(5,1).Value = 11.264
$ws2.Cells.Item(5,2).Value = 1.575
$ws2.Cells.Item(5,3).Value = 1.3
$ws2.Cells.Item(5,4).Value = 1.82

$ws2.Cells.Item(6,1).Value = 10.702
$ws2.Cells.Item(6,2).Value = 1.813
$ws2.Cells.Item(6,3).Value = 2.1
$ws2.Cells.Item(6,4).Value = 1.82

$ws2.PageSetup.Orientation = 1
$ws2.Range("A1:E6").Select()

# ----------------------------------------------------------------------
# Sheet3 "Double"
# ----------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

$ws3.Columns.Item(1).ColumnWidth = 17

$ws3.Range("F2:F7").Style = "Percent"

$ws3.PageSetup.Orientation = 1
$ws3.Range("A1:F7").Select()

# ----------------------------------------------------------------------
# Sheet4 "Regulated"
# ----------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F2:F5").Style = "Percent"

$ws4.PageSetup.Orientation = 1
$ws4.Range("A1:F5").Select()

$ws1.Activate()
